$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 140, pushing the existing rows 140-148
# down to 142-150 (same as a new weekly report being prepended to the
# top of this date-descending block).
$ws.Rows("140:141").Insert()

# --- Row 140: new "Primera" quality record ---
$ws.Range("A140").Value = 7
$ws.Range("B140").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C140").Value = "Ñuble"
$ws.Range("D140").Value = 44461
$ws.Range("D140").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E140").Value = 16
$ws.Range("F140").Value = "Fruta"
$ws.Range("G140").Value = 100101
$ws.Range("H140").Value = "Berries"
$ws.Range("I140").Value = 100101007
$ws.Range("J140").Value = "Kiwi"
$ws.Range("K140").Value = "Hayward"
$ws.Range("L140").Value = "Primera"
$ws.Range("M140").Value = 60
$ws.Range("N140").Value = 12000
$ws.Range("O140").Value = 12500
$ws.Range("P140").Value = 12250
$ws.Range("Q140").Value = "$/bandeja 18 kilos"
$ws.Range("R140").Value = "Provincia de Curicó"
$ws.Range("S140").Value = 681
$ws.Range("T140").Value = 18

# --- Row 141: new "Segunda" quality record ---
$ws.Range("A141").Value = 7
$ws.Range("B141").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C141").Value = "Ñuble"
$ws.Range("D141").Value = 44461
$ws.Range("D141").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E141").Value = 16
$ws.Range("F141").Value = "Fruta"
$ws.Range("G141").Value = 100101
$ws.Range("H141").Value = "Berries"
$ws.Range("I141").Value = 100101007
$ws.Range("J141").Value = "Kiwi"
$ws.Range("K141").Value = "Hayward"
$ws.Range("L141").Value = "Segunda"
$ws.Range("M141").Value = 60
$ws.Range("N141").Value = 11000
$ws.Range("O141").Value = 11500
$ws.Range("P141").Value = 11250
$ws.Range("Q141").Value = "$/bandeja 18 kilos"
$ws.Range("R141").Value = "Provincia de Curicó"
$ws.Range("S141").Value = 625
$ws.Range("T141").Value = 18
